$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E3").Value = 16.4482
$ws.Range("D7").Value = -7.015100000000001
$ws.Range("B10").Value = 4.708399999999996
$ws.Range("B12").Value = 5.088199999999998
$ws.Range("D15").Value = -7.904100000000002
$ws.Range("B18").Value = 7.418599999999994
$ws.Range("E18").Value = 17.6719
$ws.Range("E19").Value = 16.60850000000001
$ws.Range("D20").Value = -7.647699999999997
$ws.Range("E27").Value = 16.48599999999999
$ws.Range("D29").Value = -6.503699999999997
$ws.Range("D30").Value = -7.7709
$ws.Range("D31").Value = -8.547000000000006
$ws.Range("B37").Value = 8.938399999999996
$ws.Range("D40").Value = -8.533499999999991
$ws.Range("E42").Value = 16.61180000000001
$ws.Range("E44").Value = 16.54309999999999
$ws.Range("E47").Value = 16.3817
$ws.Range("B55").Value = 5.460399999999996
$ws.Range("E58").Value = 16.44400000000001
$ws.Range("B68").Value = 4.834199999999996
$ws.Range("D68").Value = -6.471499999999999
$ws.Range("E73").Value = 17.17860000000001
$ws.Range("D76").Value = -7.516299999999994
$ws.Range("B77").Value = 9.317700000000007
$ws.Range("B78").Value = 9.819800000000003
$ws.Range("D87").Value = -7.959199999999998
$ws.Range("D88").Value = -7.181299999999997
$ws.Range("E95").Value = 17.98560000000002
$ws.Range("D96").Value = -7.7379
$ws.Range("D98").Value = -8.726300000000004
$ws.Range("D101").Value = -7.953499999999999
$ws.Range("E101").Value = 16.6845
$ws.Range("D102").Value = -7.818799999999999
